$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(2, 2, 2, 1, 1, 1, 2, 1, 2, 2, 2, 1, 2, 2)
$cols = @("B", "C", "D", "E", "F", "G", "H", "I", "J", "K", "L", "M", "N", "O")

for ($i = 0; $i -lt $cols.Length; $i++) {
    $ws.Range("$($cols[$i])2").Value = $values[$i]
}
